$wb = $excel.ActiveWorkbook

# --- Update selection on the previously-active sheet (tryEditorCode) ---
$tryEditorCode = $wb.Worksheets.Item("tryEditorCode")
$tryEditorCode.Range("A3").Select()

# --- Add the new sheet "ArrayPracticeQnsQ1" after "tryEditorCode" ---
# Copy tryEditorCode so the new sheet inherits the workbook's normal
# sheet-level defaults (row height, etc.), then wipe its content/format.
$tryEditorCode.Copy([System.Type]::Missing, $tryEditorCode)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "ArrayPracticeQnsQ1"

# Remove all copied rows/formatting so we start from a clean sheet.
$newSheet.Range("A1:A3").EntireRow.Delete()

# Column A width
$newSheet.Columns.Item(1).ColumnWidth = 36.28

# Row 1: "pCode" label, default formatting
$newSheet.Range("A1").Value = "pCode"

# Row 2: new cleaned-up search() snippet
$code2 = "def search(input_list, num):`n    if num in input_list:`n        print(`"Element Found`")`n    else:`n        print(`"Not Found`")`nsearch([12, 23, 45, 67, 6, 90], 12)"
$newSheet.Range("A2").Value = $code2
$newSheet.Range("A2").WrapText = $true
$newSheet.Rows.Item(2).RowHeight = 100.8

# Row 3: existing (older/messier) search() snippet, reused from elsewhere in the workbook
$code3 = "def search(input_list, num):`nif(num in input_list):`nprint(`"Element Found`")`n\xc`n\xc`nelse:`nprint(`"Not Found`")`n\xc`n\xc`n\xc`n\xc`nsearch([12, 23, 45, 67, 6, 90] , 12)"
$newSheet.Range("A3").Value = $code3
$newSheet.Range("A3").WrapText = $true
$newSheet.Rows.Item(3).RowHeight = 172.8

# Selection/active cell on the new sheet
$newSheet.Range("K3").Select()
